$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.391.90"
$ws.Range("E2").Value = "  +9.91%  "
$ws.Range("D3").Value = "3.225.53"
$ws.Range("E3").Value = "  +4.16%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "398.15"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.07"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +7.02%  "
$ws.Range("E7").Value = "  +2.82%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +5.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.32"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0920"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +7.38%  "
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("D13").Value = "3.736.79"
$ws.Range("E13").Value = "  +4.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.09"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.05"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.88%  "
$ws.Range("D16").Value = "3.229.58"
$ws.Range("E16").Value = "  +4.18%  "
$ws.Range("E17").Value = "  +4.95%  "
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").Value = "56.221.62"
$ws.Range("E19").Value = "  +9.36%  "
$ws.Range("E20").Value = "  +3.19%  "
$ws.Range("E21").Value = "  +6.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "298.06"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +12.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.63"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +8.06%  "
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "28.07"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("E29").Value = "  +4.52%  "
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").Value = "  +3.59%  "
$ws.Range("E32").Value = "  +6.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0493"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "36.75"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("E35").Value = "  +2.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.39"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.11"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +25.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.52"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "134.65"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.44"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.44%  "
$ws.Range("E42").Value = "  +3.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.99"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.02%  "
$ws.Range("E44").Value = "  +3.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.284"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.22"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +55.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.18"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.129.22"
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.08"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.43"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.39%  "
$ws.Range("E51").Value = "  +11.13%  "
